$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-7) are cyclically permuted (new <- old):
# new row 2 <- old row 6
# new row 3 <- old row 7
# new row 4 <- old row 3
# new row 5 <- old row 2
# new row 6 <- old row 4
# new row 7 <- old row 5
#
# Capture the original values for columns D and K:T for rows 2-7 first,
# then write them back in the permuted order so no source data is lost
# while writing.

$sourceRows = @(6, 7, 3, 2, 4, 5)
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

$data = @{}
foreach ($r in 2..7) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $data[$r] = $rowData
}

for ($i = 0; $i -lt 6; $i++) {
    $destRow = $i + 2
    $srcRow = $sourceRows[$i]
    $srcData = $data[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcData[$col]
    }
}
